# Append the new Kaspa buy recorded on 2025-06-20 as row 34.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Column A holds the date as literal text (e.g. "06/15/2025" in the row
# above it), not an Excel date serial number. Briefly force a text number
# format so the "06/20/2025" string isn't auto-parsed into a date, then
# restore the default "Normal" style so the cell ends up unstyled just
# like the other text-date cells in this column.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "06/20/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 713.8940000000002
$ws.Cells.Item($row, 3).Value = 0.07003840906353041
$ws.Cells.Item($row, 4).Value = 50
